# Reverse the order of comma-separated "Recorded By" entries in column G
# of the active worksheet (rows 2 through the last used row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ",\s*"

        # Build the reversed list manually (index-based) since
        # [array]::Reverse() does not mutate in place in this runtime.
        $count = $parts.Count
        $reversedParts = @()
        for ($i = $count - 1; $i -ge 0; $i--) {
            $reversedParts += $parts[$i]
        }

        $newVal = [string]::Join(", ", $reversedParts)
        $cell.Value2 = $newVal
    }
}
